$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.059.28'
$ws.Range('E2').Value = '  +0.21%  '

$ws.Range('D3').Value = '1.868.54'
$ws.Range('E3').Value = '  +0.45%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.98'
$ws.Range('E5').Value = '  +0.24%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.18%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5101'
$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3876'
$ws.Range('E8').Value = '  +1.00%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08333'
$ws.Range('E9').Value = '  +1.35%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.114'
$ws.Range('E10').Value = '  +0.29%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.39'
$ws.Range('E11').Value = '  -0.32%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.199'
$ws.Range('E12').Value = '  -0.21%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.876.17'
$ws.Range('E13').Value = '  +0.75%  '

$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '20.47'
$ws.Range('E14').Value = '  -0.50%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.228'
$ws.Range('E15').Value = '  -0.39%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.36%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001096'
$ws.Range('E17').Value = '  +0.12%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.54'
$ws.Range('E18').Value = '  -0.25%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06649'
$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.68'

$ws.Range('E21').Value = '  -0.22%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.987'
$ws.Range('E22').Value = '  -0.50%  '

$ws.Range('D23').Value = '28.093.76'
$ws.Range('E23').Value = '  +0.20%  '

$ws.Range('E24').Value = '  -0.07%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.231'
$ws.Range('E25').Value = '  -0.20%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '159.07'
$ws.Range('E26').Value = '  +1.19%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.449'
$ws.Range('E27').Value = '  -2.64%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.52'
$ws.Range('E28').Value = '  +0.17%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '124.80'
$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.1051'
$ws.Range('E30').Value = '  -0.69%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.813'
$ws.Range('E32').Value = '  -1.86%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.593'
$ws.Range('E33').Value = '  -0.07%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '9.472'
$ws.Range('E34').Value = '  +0.92%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02435'
$ws.Range('E35').Value = '  +0.91%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06537'
$ws.Range('E36').Value = '  +0.30%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2208'
$ws.Range('E37').Value = '  +1.69%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.188'
$ws.Range('E38').Value = '  -0.71%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6442'
$ws.Range('E39').Value = '  -1.66%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.237'
$ws.Range('E40').Value = '  +1.22%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.964'
$ws.Range('E41').Value = '  -0.93%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.18'
$ws.Range('E42').Value = '  +0.09%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6052'
$ws.Range('E43').Value = '  -1.27%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.05'
$ws.Range('E44').Value = '  -0.50%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.672'
$ws.Range('E45').Value = '  +0.57%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.272'
$ws.Range('E46').Value = '  -0.71%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.004'
$ws.Range('E47').Value = '  -0.45%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.227'
$ws.Range('E48').Value = '  +1.64%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '120.47'
$ws.Range('E49').Value = '  +0.27%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06887'
$ws.Range('E50').Value = '  +0.73%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '77.66'
$ws.Range('E51').Value = '  -0.92%  '
